$wb = $excel.ActiveWorkbook

# The "Users" sheet holds test-user rows; row 2 maps a CAO-type user to a
# display name. Update the display name used for the CaoUser test account.
$usersSheet = $wb.Worksheets.Item("Users")
$usersSheet.Range("A2").Value = "Drew Koecher"

# Reflect the author's final view state: they ended up on the Users sheet
# with cell B10 selected.
$usersSheet.Activate()
$usersSheet.Range("B10").Select()
